# Edit script: add storage_content__Building1 sheet, fix labels/values in
# env_impacts and capStorages sheets, update a couple of flow values in
# electricityBus__Building1, and drop the storage_content column from
# shSourceBus__Building1 (moved into its own new sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. electricityBus__Building1: update B13/C13/B23/C23
# ---------------------------------------------------------------------
$wsElectricityBus = $wb.Worksheets.Item("electricityBus__Building1")
$wsElectricityBus.Range("B13").Value = 15.0274582141794
$wsElectricityBus.Range("C13").Value = 61.5423219105333
$wsElectricityBus.Range("B23").Value = 10.44800906919004
$wsElectricityBus.Range("C23").Value = 0

# ---------------------------------------------------------------------
# 2. Create the new storage_content__Building1 sheet by duplicating the
#    shSourceBus__Building1 sheet (same date column + formatting), then
#    drop the extra columns and fill in the storage-content values.
# ---------------------------------------------------------------------
$wsShSourceBus = $wb.Worksheets.Item("shSourceBus__Building1")
$wsCosts = $wb.Worksheets.Item("costs__Building1")
$wsShSourceBus.Copy($wsCosts)
$wsStorageContent = $wb.Worksheets.Item("shSourceBus__Building1 (2)")
$wsStorageContent.Name = "storage_content__Building1"

# Drop columns B,C,D,E (flow columns) but keep column A (the date axis),
# then rebuild B (shStorage) and C (electricalStorage) with the new data.
$wsStorageContent.Range("D1:E26").Delete()
$wsStorageContent.Range("B1").Value = "shStorage__B001_storage_content"
$wsStorageContent.Range("C1").Value = "electricalStorage__B001_storage_content"

$storageContentRows = @(
    @{ Row = 2;  B = 24.41023943949497;  C = 109.3557807313954 },
    @{ Row = 3;  B = 11.13843355808191;  C = 108.8097009976744 },
    @{ Row = 4;  B = 26.64435025980561;  C = 100.3848437329652 },
    @{ Row = 5;  B = 13.36468957599664;  C = 99.74579057715127 },
    @{ Row = 6;  B = 0;                  C = 99.01706964691871 },
    @{ Row = 7;  B = 0;                  C = 94.63202081251046 },
    @{ Row = 8;  B = 0;                  C = 87.41951374037683 },
    @{ Row = 9;  B = 0;                  C = 80.9903867504638 },
    @{ Row = 10; B = 0;                  C = 74.12152314320944 },
    @{ Row = 11; B = 0;                  C = 70.4654556743557 },
    @{ Row = 12; B = 0;                  C = 70.4654556743557 },
    @{ Row = 13; B = 0;                  C = 52.99166705321686 },
    @{ Row = 14; B = 52.24448642575736;  C = 52.99166705321686 },
    @{ Row = 15; B = 49.10231151640725;  C = 52.99166705321686 },
    @{ Row = 16; B = 52.325;             C = 49.70673626925328 },
    @{ Row = 17; B = 52.325;             C = 45.97127300138378 },
    @{ Row = 18; B = 44.3433705550596;   C = 44.34516004440703 },
    @{ Row = 19; B = 35.31286976631884;  C = 42.20638928045354 },
    @{ Row = 20; B = 25.83785739473363;  C = 39.87401386533726 },
    @{ Row = 21; B = 18.49078851382524;  C = 29.81635605836052 },
    @{ Row = 22; B = 11.81135031878672;  C = 15.73085772115121 },
    @{ Row = 23; B = 32.39808334314499;  C = 3.58200996627907 },
    @{ Row = 24; B = 22.9580788963761;   C = 1.617325581395349 },
    @{ Row = 25; B = 11.94378045273504;  C = 0.6711960127906977 },
    @{ Row = 26; B = 0;                  C = 0 }
)

foreach ($entry in $storageContentRows) {
    $r = $entry.Row
    $wsStorageContent.Range("B$r").Value = $entry.B
    $wsStorageContent.Range("C$r").Value = $entry.C
}

# ---------------------------------------------------------------------
# 3. shSourceBus__Building1: remove the storage_content column (E) now
#    that it lives in its own sheet.
# ---------------------------------------------------------------------
$wsShSourceBus.Range("E1:E26").Delete()

# ---------------------------------------------------------------------
# 4. costs__Building1: tiny floating point correction to Feed-in value.
# ---------------------------------------------------------------------
$wsCosts.Range("B5").Value = -7.034695744793923

# ---------------------------------------------------------------------
# 5. env_impacts__Building1: swap the electricalStorage/shStorage rows.
# ---------------------------------------------------------------------
$wsEnvImpacts = $wb.Worksheets.Item("env_impacts__Building1")
$wsEnvImpacts.Range("A6").Value = "shStorage__Building1"
$wsEnvImpacts.Range("B6").Value = 22.37625
$wsEnvImpacts.Range("A7").Value = "electricalStorage__Building1"
$wsEnvImpacts.Range("B7").Value = 2757.066666666667

# ---------------------------------------------------------------------
# 6. capStorages__Building1: swap the electricalStorage/shStorage rows.
# ---------------------------------------------------------------------
$wsCapStorages = $wb.Worksheets.Item("capStorages__Building1")
$wsCapStorages.Range("A2").Value = "shStorage__Building1"
$wsCapStorages.Range("B2").Value = 4500
$wsCapStorages.Range("A3").Value = "electricalStorage__Building1"
$wsCapStorages.Range("B3").Value = 200

$wb.Save()
